# Update data: 2025-10-29 16:03
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 04:03 PM"

# --- Top Losers sheet: a new row (TBOTEK) enters the top-losers list at
#     row 60, pushing the existing rows 60-66 down by one, and the prior
#     TBOTEK entry (old row 67) drops out of this window. ---
$losers = $wb.Worksheets.Item("Top Losers")

# Snapshot rows 60-66 (7 rows) before overwriting them, since row 67 needs
# the values currently sitting in row 66, row 66 needs row 65's, etc.
$startRow = 60
$endRow = 66
$rowCount = $endRow - $startRow + 1

$icon = @()
$stock = @()
$latest = @()
$weekly = @()
$monthly = @()

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $icon += $losers.Cells.Item($r, 1).Value2
    $stock += $losers.Cells.Item($r, 2).Value2
    $latest += $losers.Cells.Item($r, 3).Value2
    $weekly += $losers.Cells.Item($r, 4).Value2
    $monthly += $losers.Cells.Item($r, 5).Value2
}

# Shift the captured rows down by one (row 61 gets old row 60's data, ...,
# row 67 gets old row 66's data).
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i + 1
    $losers.Cells.Item($r, 1).Value = $icon[$i]
    $losers.Cells.Item($r, 2).Value = $stock[$i]
    $losers.Cells.Item($r, 3).Value = $latest[$i]
    $losers.Cells.Item($r, 4).Value = $weekly[$i]
    $losers.Cells.Item($r, 5).Value = $monthly[$i]
}

# Write the new entry into row 60 (reuse the existing "down" icon glyph
# captured above rather than retyping the emoji literal - every row in
# this sheet uses the same icon).
$losers.Cells.Item($startRow, 1).Value = $icon[0]
$losers.Cells.Item($startRow, 2).Value = "TBOTEK"
$losers.Cells.Item($startRow, 3).Value = -2.524
$losers.Cells.Item($startRow, 4).Value = -3.5732
$losers.Cells.Item($startRow, 5).Value = 1.036
